$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the numeric id values in column A with composite "id|country|place|date"
# strings. Rows 3-19 are written first, then row 2, so new shared-string entries
# land in the same order as in the target workbook.
$values = [ordered]@{
    3  = "105671|SE|Haessleholm|2016-12-01"
    4  = "79179|SE|Filipstad|2017-01-01"
    5  = "79424|SE|Kalmar|2017-03-01"
    6  = "79612|SE|Joenkoeping|2017-05-01"
    7  = "83364|SE|Lidkoeping|2017-07-01"
    8  = "83367|SE|Lidkoeping|2017-07-01"
    9  = "83376|SE|Oerebro|2017-08-01"
    10 = "70319|SE|Oerebro|2018-01-01"
    11 = "70320|SE|Malmo|2018-06-01"
    12 = "79180|SE|Filipstad|2018-07-01"
    13 = "79614|SE|Kalmar|2018-08-01"
    14 = "79615|SE|Kalmar|2018-08-01"
    15 = "107204|SE|Oerebro|2018-08-01"
    16 = "110224|SE|Oerebro|2018-08-01"
    17 = "105668|SE|Joenkoeping|2019-07-01"
    18 = "105670|SE|Malmoe|2019-07-01"
    19 = "110223|SE|Joenkoeping|2019-08-01"
    2  = "79418|SE|Hoeoer|2016-09-01"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}

# Widen column A to fit the new longer strings (target stored width = 40;
# the persisted XML width adds a fixed ~0.8333 padding on top of the
# ColumnWidth value actually applied, so compensate here).
$ws.Columns.Item(1).ColumnWidth = 39.166666666666664

# Update the selected cell.
$ws.Range("A12").Select()
